$d = $word.ActiveDocument

$d.Content.Find.Execute("71×28=1988", $true, $false, $false, $false, $false, $true, 1, $false, "17×57=969", 2) | Out-Null
$d.Content.Find.Execute("47×38=1786", $true, $false, $false, $false, $false, $true, 1, $false, "70×11=770", 2) | Out-Null
$d.Content.Find.Execute("68×86=5848", $true, $false, $false, $false, $false, $true, 1, $false, "86×91=7826", 2) | Out-Null
$d.Content.Find.Execute("73×66=4818", $true, $false, $false, $false, $false, $true, 1, $false, "13×68=884", 2) | Out-Null
$d.Content.Find.Execute("50×47=2350", $true, $false, $false, $false, $false, $true, 1, $false, "79×28=2212", 2) | Out-Null
$d.Content.Find.Execute("33×35=1155", $true, $false, $false, $false, $false, $true, 1, $false, "92×49=4508", 2) | Out-Null
$d.Content.Find.Execute("62×38=2356", $true, $false, $false, $false, $false, $true, 1, $false, "34×33=1122", 2) | Out-Null
$d.Content.Find.Execute("63×78=4914", $true, $false, $false, $false, $false, $true, 1, $false, "74×48=3552", 2) | Out-Null
$d.Content.Find.Execute("95×60=5700", $true, $false, $false, $false, $false, $true, 1, $false, "32×17=544", 2) | Out-Null
$d.Content.Find.Execute("59×28=1652", $true, $false, $false, $false, $false, $true, 1, $false, "50×22=1100", 2) | Out-Null
$d.Content.Find.Execute("57×43=2451", $true, $false, $false, $false, $false, $true, 1, $false, "13×46=598", 2) | Out-Null
$d.Content.Find.Execute("96×21=2016", $true, $false, $false, $false, $false, $true, 1, $false, "91×35=3185", 2) | Out-Null
$d.Content.Find.Execute("59×16=944", $true, $false, $false, $false, $false, $true, 1, $false, "89×12=1068", 2) | Out-Null
$d.Content.Find.Execute("89×80=7120", $true, $false, $false, $false, $false, $true, 1, $false, "37×89=3293", 2) | Out-Null
$d.Content.Find.Execute("55×77=4235", $true, $false, $false, $false, $false, $true, 1, $false, "12×54=648", 2) | Out-Null
$d.Content.Find.Execute("82×59=4838", $true, $false, $false, $false, $false, $true, 1, $false, "94×32=3008", 2) | Out-Null
$d.Content.Find.Execute("58×63=3654", $true, $false, $false, $false, $false, $true, 1, $false, "62×56=3472", 2) | Out-Null
$d.Content.Find.Execute("38×85=3230", $true, $false, $false, $false, $false, $true, 1, $false, "83×89=7387", 2) | Out-Null
$d.Content.Find.Execute("95×59=5605", $true, $false, $false, $false, $false, $true, 1, $false, "91×68=6188", 2) | Out-Null
$d.Content.Find.Execute("96×97=9312", $true, $false, $false, $false, $false, $true, 1, $false, "73×14=1022", 2) | Out-Null
$d.Content.Find.Execute("32×51=1632", $true, $false, $false, $false, $false, $true, 1, $false, "40×79=3160", 2) | Out-Null
$d.Content.Find.Execute("15×45=675", $true, $false, $false, $false, $false, $true, 1, $false, "86×11=946", 2) | Out-Null
$d.Content.Find.Execute("28×80=2240", $true, $false, $false, $false, $false, $true, 1, $false, "69×67=4623", 2) | Out-Null
$d.Content.Find.Execute("77×36=2772", $true, $false, $false, $false, $false, $true, 1, $false, "89×89=7921", 2) | Out-Null
$d.Content.Find.Execute("14×22=308", $true, $false, $false, $false, $false, $true, 1, $false, "28×98=2744", 2) | Out-Null
